# Bug fix for one trial with wrong date
# Updates summary_within_12m_after_CD counts/percentages/CIs for several
# institutions (and the Total row) in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Aalborg University Hospital
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 11.4
$ws.Range("E3").Value = 4.5
$ws.Range("F3").Value = 26

# Row 4: Aarhus University
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 2.1
$ws.Range("E4").Value = 0.8
$ws.Range("F4").Value = 5.3

# Row 5: Aarhus University Hospital
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 6.1
$ws.Range("E5").Value = 2.4
$ws.Range("F5").Value = 14.6

# Row 7: Bispebjerg and Frederiksberg Hospital
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 3.9
$ws.Range("E7").Value = 1.3
$ws.Range("F7").Value = 10.8

# Row 32: Skane University Hospital
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 4.3
$ws.Range("E32").Value = 0.2
$ws.Range("F32").Value = 21

# Row 56: Total
$ws.Range("C56").Value = 48
$ws.Range("D56").Value = 2.3
$ws.Range("E56").Value = 1.7
$ws.Range("F56").Value = 3
